# eprepago.xlsx - configure data per commit:
# "se configura la data en el excel eprepago"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Remove the two mailto hyperlinks (N2, N3) from the sheet.
$ws.Hyperlinks.Delete()

# Drop the correo / tipoCorreo / numeroCelular columns (N, O, P) of the
# header row and the numeroCelular value column (P) in the data rows,
# while keeping the (now value-less) styled cells N2:O3.
$ws.Range("N1:P1").ClearContents()
$ws.Range("N2:O3").ClearContents()
$ws.Range("P2:P3").ClearContents()

# D2 already shows "autotest27" (its shared string just gets renumbered
# automatically); make D3 match so that the old "userunico01" string is
# no longer referenced anywhere and drops out of the shared string table.
$ws.Range("D3").Value = "autotest27"

# Update the document number in row 3 and copy B2's number formatting.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = 93221452
$excel.CutCopyMode = 0

# Update the view: scroll right and move the selection to S3.
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("S3").Select() | Out-Null
